# Insert a new price-record row (row 168) into the "Naranja" sheet, which
# shifts the existing rows 168-181 down to 169-182 (the trailing row 182
# ends up holding what used to be row 181's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 168..181 down by one row, inheriting row 168's formatting
# (this is what gives the new row its date-formatted column D style).
$ws.Rows("168:168").Insert()

# Populate the newly-opened row 168 with the new record.
$ws.Cells.Item(168, 1).Value  = 1
$ws.Cells.Item(168, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(168, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(168, 4).Value  = 45194
$ws.Cells.Item(168, 5).Value  = 15
$ws.Cells.Item(168, 6).Value  = "Fruta"
$ws.Cells.Item(168, 7).Value  = 100102
$ws.Cells.Item(168, 8).Value  = "Cítricos"
$ws.Cells.Item(168, 9).Value  = 100102005
$ws.Cells.Item(168, 10).Value = "Naranja"
$ws.Cells.Item(168, 11).Value = "Navel Late"
$ws.Cells.Item(168, 12).Value = "Segunda"
$ws.Cells.Item(168, 13).Value = 300
$ws.Cells.Item(168, 14).Value = 750
$ws.Cells.Item(168, 15).Value = 800
$ws.Cells.Item(168, 16).Value = 775
$ws.Cells.Item(168, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(168, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 19).Value = 775
$ws.Cells.Item(168, 20).Value = 1
